$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Viking")
$ws2 = $wb.Worksheets.Item("NAF")

# --- Free up the two shared strings that are no longer referenced anywhere
# ("Location" header and its Google-maps URL, previously NAF!E1 / NAF!E7) so
# the new Longitude/Latitude data lands in the same slots the original file
# used (and the orphaned strings drop out of the shared-string table).
$ws2.Range("E1").ClearContents()
$ws2.Range("E7").ClearContents()

# --- Viking sheet: add Longitude / Latitude columns (Vollaveien / Oslo) ---
$ws1.Range("E2").NumberFormat = "@"
$ws1.Range("E2").Value = "59.92815"
$ws1.Range("E2").ClearFormats()

$ws1.Range("F2").NumberFormat = "@"
$ws1.Range("F2").Value = "10.84207"
$ws1.Range("F2").ClearFormats()

$ws1.Range("E1").Value = "Longitude"
$ws1.Range("F1").Value = "Latitude"

# --- NAF sheet: add Longitude / Latitude values for Trondheim and Sandnes ---
$ws2.Range("E2").NumberFormat = "@"
$ws2.Range("E2").Value = "63.35462"
$ws2.Range("E2").ClearFormats()

$ws2.Range("F2").NumberFormat = "@"
$ws2.Range("F2").Value = "10.37234"
$ws2.Range("F2").ClearFormats()

$ws2.Range("E3").NumberFormat = "@"
$ws2.Range("E3").Value = "58.87855"
$ws2.Range("E3").ClearFormats()

$ws2.Range("F3").NumberFormat = "@"
$ws2.Range("F3").Value = "5.71927"
$ws2.Range("F3").ClearFormats()

# Re-create the NAF header row's Longitude/Latitude labels (re-uses the
# shared strings created above for Viking's headers).
$ws2.Range("E1").Value = "Longitude"
$ws2.Range("F1").Value = "Latitude"

# --- View state: selection / zoom to match the saved workbook ---
$ws1.Select()
$ws1.Range("E1:F1").Select()
$excel.ActiveWindow.Zoom = 140

$ws2.Select()
$ws2.Range("E4").Select()
